$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$ws2 = $wb.Worksheets.Item("Summary")

# Append new order rows (32-41) in column C (PackageID/FlowerName text column)
$ws.Range("C32").Value = "100_绣球单瓣白_Hydrangea White S_Hydrangea L._1stem"
$ws.Range("C33").Value = "107_绣球单瓣浅粉_Hydrangea Light Pink S_Hydrangea L._1stem"
$ws.Range("C34").Value = "225_果汁阳台_Juicy Terrazza_Rosa rugosa Thunb._10stems"
$ws.Range("C35").Value = "277_草莓杏仁饼_undefined_Rosa rugosa Thunb._10stems"
$ws.Range("C36").Value = "238_苏菲宝贝_undefined_Rosa rugosa Thunb._10stems"
$ws.Range("C37").Value = "244_繁星_undefined_Rosa rugosa Thunb._10stems"
$ws.Range("C38").Value = "274_仙子之吻_undefined_Rosa rugosa Thunb._10stems"
$ws.Range("C39").Value = "276_情迷罗拉_undefined_Rosa rugosa Thunb._10stems"
$ws.Range("C40").Value = "221_朱丽叶塔_Julieta_Rosa rugosa Thunb._10stems"
$ws.Range("C41").Value = "597_尤加利叶小叶_undefined_undefined_1bunch"

# Update the Summary sheet's concatenated TotalNumber tracking string (G2).
# The new value is all-digits, which Excel would otherwise auto-convert to a
# number. To preserve it as text (as in the source file, t="str") without
# leaving a quote-prefix style behind on the destination cell, stage the text
# (with a leading apostrophe to force text entry) in a scratch cell that sits
# well within the sheet's final used range, copy it, and paste-special just
# the values into G2 - this keeps G2's type as text but carries no style.
# The scratch cell is fully cleared afterwards so it leaves no trace.
$ws.Range("L2").Value = "'052000000000000000000000000000000000000000"
$ws.Range("L2").Copy()
$ws2.Range("G2").PasteSpecial(-4163)
$ws.Range("L2").Clear()
